# Apply updated cryptocurrency price/volume data to Sheet1 (rows 2-51)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.515.77"
$ws.Range("E2").Value = "  +3.63%  "

$ws.Range("D3").Value = "1.587.09"
$ws.Range("E3").Value = "  +0.86%  "

$ws.Range("E4").Value = "  +1.01%  "

$ws.Range("D5").Value = "'212.97"
$ws.Range("E5").Value = "  +0.66%  "

$ws.Range("D6").Value = "'0.490"
$ws.Range("E6").Value = "  -0.54%  "

$ws.Range("E7").Value = "  +1.01%  "

$ws.Range("D8").Value = "'24.27"
$ws.Range("E8").Value = "  +6.11%  "

$ws.Range("D9").Value = "'0.252"
$ws.Range("E9").Value = "  +0.41%  "

$ws.Range("E10").Value = "  +0.17%  "

$ws.Range("E11").Value = "  +1.75%  "

$ws.Range("D12").Value = "1.815.14"
$ws.Range("E12").Value = "  +1.06%  "

$ws.Range("D13").Value = "1.596.04"
$ws.Range("E13").Value = "  +1.58%  "

$ws.Range("D14").Value = "'0.529"
$ws.Range("E14").Value = "  +1.71%  "

$ws.Range("E15").Value = "  -0.57%  "

$ws.Range("D16").Value = "28.508.62"
$ws.Range("E16").Value = "  +3.71%  "

$ws.Range("E17").Value = "  +0.98%  "

$ws.Range("D18").Value = "'231.57"
$ws.Range("E18").Value = "  +2.05%  "

$ws.Range("E19").Value = "  -0.08%  "

$ws.Range("E20").Value = "  -0.81%  "

$ws.Range("E21").Value = "  +0.97%  "

$ws.Range("D22").Value = "'4.04"
$ws.Range("E22").Value = "  -2.14%  "

$ws.Range("E24").Value = "  +1.68%  "

$ws.Range("D25").Value = "'151.75"
$ws.Range("E25").Value = "  +0.66%  "

$ws.Range("D26").Value = "'15.24"
$ws.Range("E26").Value = "  +0.37%  "

$ws.Range("E27").Value = "  -0.71%  "

$ws.Range("E28").Value = "  -0.88%  "

$ws.Range("E29").Value = "  +0.92%  "

$ws.Range("E30").Value = "  -1.33%  "

$ws.Range("E31").Value = "  -0.38%  "

$ws.Range("E32").Value = "  -0.03%  "

$ws.Range("E33").Value = "  +0.73%  "

$ws.Range("D34").Value = "1.400.84"
$ws.Range("E34").Value = "  -3.88%  "

$ws.Range("E35").Value = "  -1.32%  "

$ws.Range("D36").Value = "'1.02"
$ws.Range("E36").Value = "  -9.94%  "

$ws.Range("E37").Value = "  +1.06%  "

$ws.Range("E38").Value = "  +10.57%  "

$ws.Range("E39").Value = "  -0.55%  "

$ws.Range("E40").Value = "  +0.32%  "

$ws.Range("D41").Value = "'0.813"
$ws.Range("E41").Value = "  -0.09%  "

$ws.Range("E42").Value = "  +0.98%  "

$ws.Range("E43").Value = "  -0.31%  "

$ws.Range("D44").Value = "'1.86"
$ws.Range("E44").Value = "  +0.37%  "

$ws.Range("D45").Value = "'0.983"
$ws.Range("E45").Value = "  +0.48%  "

$ws.Range("D46").Value = "'63.00"
$ws.Range("E46").Value = "  -2.07%  "

$ws.Range("D47").Value = "1.725.47"
$ws.Range("E47").Value = "  +1.05%  "

$ws.Range("E48").Value = "  +1.55%  "

$ws.Range("D49").Value = "'87.29"
$ws.Range("E49").Value = "  +0.34%  "

$ws.Range("D50").Value = "0.0₆0105"
$ws.Range("E50").Value = "  +0.95%  "

$ws.Range("E51").Value = "  -0.67%  "
